$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 44.30160133333333
$ws.Range("H2").Value = 132.904804
$ws.Range("I2").Value = 0.1310981746002036
$ws.Range("J2").Value = 0.1395903267786693
$ws.Range("M2").Value = 0.09777333333333334
$ws.Range("N2").Value = 0.29332
$ws.Range("O2").Value = 0.3178747086843887
$ws.Range("P2").Value = 0.3620126183588235
$ws.Range("Q2").Value = 4.331515234364445
$ws.Range("R2").Value = 38.98363710928
$ws.Range("S2").Value = 0.04167279406009484
$ws.Range("T2").Value = 0.05053345969470988
$ws.Range("G3").Value = 44.30160133333333
$ws.Range("H3").Value = 132.904804
$ws.Range("I3").Value = 0.1310981746002036
$ws.Range("J3").Value = 0.1395903267786693
$ws.Range("O3").Value = 0.1339783593343184
$ws.Range("P3").Value = 0.1525816786959054
$ws.Range("Q3").Value = 1.825654223746222
$ws.Range("R3").Value = 16.430888013716
$ws.Range("S3").Value = 0.0175643183446593
$ws.Range("T3").Value = 0.02129892638959937
$ws.Range("G4").Value = 44.30160133333333
$ws.Range("H4").Value = 132.904804
$ws.Range("I4").Value = 0.1310981746002036
$ws.Range("J4").Value = 0.1395903267786693
$ws.Range("M4").Value = 0.025476
$ws.Range("N4").Value = 0.076428
$ws.Range("O4").Value = 0.08282602016681595
$ws.Range("P4").Value = 0.09432667528954096
$ws.Range("Q4").Value = 1.128627595568
$ws.Range("R4").Value = 10.157648360112
$ws.Range("S4").Value = 0.01085834005326922
$ws.Range("T4").Value = 0.01316709142761246
$ws.Range("G5").Value = 44.30160133333333
$ws.Range("H5").Value = 132.904804
$ws.Range("I5").Value = 0.1310981746002036
$ws.Range("J5").Value = 0.1395903267786693
$ws.Range("M5").Value = 0.1125055
$ws.Range("N5").Value = 0.225011
$ws.Range("O5").Value = 0.3657710320253459
$ws.Range("P5").Value = 0.2777063318885082
$ws.Range("Q5").Value = 4.984173808807333
$ws.Range("R5").Value = 29.905042852844
$ws.Range("S5").Value = 0.04795191462015547
$ws.Range("T5").Value = 0.03876511761682246
$ws.Range("G6").Value = 44.30160133333333
$ws.Range("H6").Value = 132.904804
$ws.Range("I6").Value = 0.1310981746002036
$ws.Range("J6").Value = 0.1395903267786693
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.03062
$ws.Range("N6").Value = 0.09186
$ws.Range("O6").Value = 0.09954987978913112
$ws.Range("P6").Value = 0.1133726957672219
$ws.Range("Q6").Value = 1.356515032826666
$ws.Range("R6").Value = 12.20863529544
$ws.Range("S6").Value = 0.01305080752202479
$ws.Range("T6").Value = 0.01582573164992516
$ws.Range("I7").Value = 0.2142454163706631
$ws.Range("J7").Value = 0.2281236010586413
$ws.Range("M7").Value = 0.09777333333333334
$ws.Range("N7").Value = 0.29332
$ws.Range("O7").Value = 0.3178747086843887
$ws.Range("P7").Value = 0.3620126183588235
$ws.Range("Q7").Value = 7.078720109813334
$ws.Range("R7").Value = 63.70848098832001
$ws.Range("S7").Value = 0.06810319931579009
$ws.Range("T7").Value = 0.08258362212868241
$ws.Range("I8").Value = 0.2142454163706631
$ws.Range("J8").Value = 0.2281236010586413
$ws.Range("O8").Value = 0.1339783593343184
$ws.Range("P8").Value = 0.1525816786959054
$ws.Range("S8").Value = 0.02870424938023937
$ws.Range("T8").Value = 0.03480748199968252
$ws.Range("I9").Value = 0.2142454163706631
$ws.Range("J9").Value = 0.2281236010586413
$ws.Range("M9").Value = 0.025476
$ws.Range("N9").Value = 0.076428
$ws.Range("O9").Value = 0.08282602016681595
$ws.Range("P9").Value = 0.09432667528954096
$ws.Range("Q9").Value = 1.844444362992
$ws.Range("R9").Value = 16.599999266928
$ws.Range("S9").Value = 0.01774509517696442
$ws.Range("T9").Value = 0.02151814084293924
$ws.Range("I10").Value = 0.2142454163706631
$ws.Range("J10").Value = 0.2281236010586413
$ws.Range("M10").Value = 0.1125055
$ws.Range("N10").Value = 0.225011
$ws.Range("O10").Value = 0.3657710320253459
$ws.Range("P10").Value = 0.2777063318885082
$ws.Range("Q10").Value = 8.145318546106001
$ws.Range("R10").Value = 48.87191127663601
$ws.Range("S10").Value = 0.0783647670525974
$ws.Range("T10").Value = 0.06335136846719266
$ws.Range("I11").Value = 0.2142454163706631
$ws.Range("J11").Value = 0.2281236010586413
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.03062
$ws.Range("N11").Value = 0.09186
$ws.Range("O11").Value = 0.09954987978913112
$ws.Range("P11").Value = 0.1133726957672219
$ws.Range("Q11").Value = 2.21686632104
$ws.Range("R11").Value = 19.95179688936
$ws.Range("S11").Value = 0.02132810544507186
$ws.Range("T11").Value = 0.02586298762014443
$ws.Range("G12").Value = 82.35175066666666
$ws.Range("H12").Value = 247.055252
$ws.Range("I12").Value = 0.2436969288378267
$ws.Range("J12").Value = 0.2594828954344383
$ws.Range("M12").Value = 0.09777333333333334
$ws.Range("N12").Value = 0.29332
$ws.Range("O12").Value = 0.3178747086843887
$ws.Range("P12").Value = 0.3620126183588235
$ws.Range("Q12").Value = 8.051805168515555
$ws.Range("R12").Value = 72.46624651664
$ws.Range("S12").Value = 0.07746509026160435
$ws.Range("T12").Value = 0.09393608239554979
$ws.Range("G13").Value = 82.35175066666666
$ws.Range("H13").Value = 247.055252
$ws.Range("I13").Value = 0.2436969288378267
$ws.Range("J13").Value = 0.2594828954344383
$ws.Range("O13").Value = 0.1339783593343184
$ws.Range("P13").Value = 0.1525816786959054
$ws.Range("Q13").Value = 3.393688194389777
$ws.Range("R13").Value = 30.543193749508
$ws.Range("S13").Value = 0.03265011470050417
$ws.Range("T13").Value = 0.03959233577826068
$ws.Range("G14").Value = 82.35175066666666
$ws.Range("H14").Value = 247.055252
$ws.Range("I14").Value = 0.2436969288378267
$ws.Range("J14").Value = 0.2594828954344383
$ws.Range("M14").Value = 0.025476
$ws.Range("N14").Value = 0.076428
$ws.Range("O14").Value = 0.08282602016681595
$ws.Range("P14").Value = 0.09432667528954096
$ws.Range("Q14").Value = 2.097993199984
$ws.Range("R14").Value = 18.881938799856
$ws.Range("S14").Value = 0.02018444674251295
$ws.Range("T14").Value = 0.02447615882083417
$ws.Range("G15").Value = 82.35175066666666
$ws.Range("H15").Value = 247.055252
$ws.Range("I15").Value = 0.2436969288378267
$ws.Range("J15").Value = 0.2594828954344383
$ws.Range("M15").Value = 0.1125055
$ws.Range("N15").Value = 0.225011
$ws.Range("O15").Value = 0.3657710320253459
$ws.Range("P15").Value = 0.2777063318885082
$ws.Range("Q15").Value = 9.265024884628666
$ws.Range("R15").Value = 55.590149307772
$ws.Range("S15").Value = 0.08913727716241915
$ws.Range("T15").Value = 0.07206004307890718
$ws.Range("G16").Value = 82.35175066666666
$ws.Range("H16").Value = 247.055252
$ws.Range("I16").Value = 0.2436969288378267
$ws.Range("J16").Value = 0.2594828954344383
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03062
$ws.Range("N16").Value = 0.09186
$ws.Range("O16").Value = 0.09954987978913112
$ws.Range("P16").Value = 0.1133726957672219
$ws.Range("Q16").Value = 2.521610605413333
$ws.Range("R16").Value = 22.69449544872
$ws.Range("S16").Value = 0.02425999997078608
$ws.Range("T16").Value = 0.02941827536088642
$ws.Range("G17").Value = 61.6746195
$ws.Range("H17").Value = 123.349239
$ws.Range("I17").Value = 0.1825087534596294
$ws.Range("J17").Value = 0.1295540872992837
$ws.Range("M17").Value = 0.09777333333333334
$ws.Range("N17").Value = 0.29332
$ws.Range("O17").Value = 0.3178747086843887
$ws.Range("P17").Value = 0.3620126183588235
$ws.Range("Q17").Value = 6.03013313058
$ws.Range("R17").Value = 36.18079878348
$ws.Range("S17").Value = 0.05801491683833061
$ws.Range("T17").Value = 0.04690021436230129
$ws.Range("G18").Value = 61.6746195
$ws.Range("H18").Value = 123.349239
$ws.Range("I18").Value = 0.1825087534596294
$ws.Range("J18").Value = 0.1295540872992837
$ws.Range("O18").Value = 0.1339783593343184
$ws.Range("P18").Value = 0.1525816786959054
$ws.Range("Q18").Value = 2.5415905113885
$ws.Range("R18").Value = 15.249543068331
$ws.Range("S18").Value = 0.02445222335267276
$ws.Range("T18").Value = 0.01976758012204059
$ws.Range("G19").Value = 61.6746195
$ws.Range("H19").Value = 123.349239
$ws.Range("I19").Value = 0.1825087534596294
$ws.Range("J19").Value = 0.1295540872992837
$ws.Range("M19").Value = 0.025476
$ws.Range("N19").Value = 0.076428
$ws.Range("O19").Value = 0.08282602016681595
$ws.Range("P19").Value = 0.09432667528954096
$ws.Range("Q19").Value = 1.571222606382
$ws.Range("R19").Value = 9.427335638292
$ws.Range("S19").Value = 0.01511647369466771
$ws.Range("T19").Value = 0.01222040632511238
$ws.Range("G20").Value = 61.6746195
$ws.Range("H20").Value = 123.349239
$ws.Range("I20").Value = 0.1825087534596294
$ws.Range("J20").Value = 0.1295540872992837
$ws.Range("M20").Value = 0.1125055
$ws.Range("N20").Value = 0.225011
$ws.Range("O20").Value = 0.3657710320253459
$ws.Range("P20").Value = 0.2777063318885082
$ws.Range("Q20").Value = 6.93873390415725
$ws.Range("R20").Value = 27.754935616629
$ws.Range("S20").Value = 0.06675641510658807
$ws.Range("T20").Value = 0.03597799036504764
$ws.Range("G21").Value = 61.6746195
$ws.Range("H21").Value = 123.349239
$ws.Range("I21").Value = 0.1825087534596294
$ws.Range("J21").Value = 0.1295540872992837
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.03062
$ws.Range("N21").Value = 0.09186
$ws.Range("O21").Value = 0.09954987978913112
$ws.Range("P21").Value = 0.1133726957672219
$ws.Range("Q21").Value = 1.88847684909
$ws.Range("R21").Value = 11.33086109454
$ws.Range("S21").Value = 0.01816872446737028
$ws.Range("T21").Value = 0.01468789612478179
$ws.Range("G22").Value = 77.19964866666666
$ws.Range("H22").Value = 231.598946
$ws.Range("I22").Value = 0.2284507267316773
$ws.Range("J22").Value = 0.2432490894289675
$ws.Range("M22").Value = 0.09777333333333334
$ws.Range("N22").Value = 0.29332
$ws.Range("O22").Value = 0.3178747086843887
$ws.Range("P22").Value = 0.3620126183588235
$ws.Range("Q22").Value = 7.548066982302222
$ws.Range("R22").Value = 67.93260284071999
$ws.Range("S22").Value = 0.07261870820856879
$ws.Range("T22").Value = 0.08805923977758014
$ws.Range("G23").Value = 77.19964866666666
$ws.Range("H23").Value = 231.598946
$ws.Range("I23").Value = 0.2284507267316773
$ws.Range("J23").Value = 0.2432490894289675
$ws.Range("O23").Value = 0.1339783593343184
$ws.Range("P23").Value = 0.1525816786959054
$ws.Range("Q23").Value = 3.181371788337111
$ws.Range("R23").Value = 28.632346095034
$ws.Range("S23").Value = 0.03060745355624284
$ws.Range("T23").Value = 0.03711535440632229
$ws.Range("G24").Value = 77.19964866666666
$ws.Range("H24").Value = 231.598946
$ws.Range("I24").Value = 0.2284507267316773
$ws.Range("J24").Value = 0.2432490894289675
$ws.Range("M24").Value = 0.025476
$ws.Range("N24").Value = 0.076428
$ws.Range("O24").Value = 0.08282602016681595
$ws.Range("P24").Value = 0.09432667528954096
$ws.Range("Q24").Value = 1.966738249432
$ws.Range("R24").Value = 17.700644244888
$ws.Range("S24").Value = 0.01892166449940166
$ws.Range("T24").Value = 0.02294487787304273
$ws.Range("G25").Value = 77.19964866666666
$ws.Range("H25").Value = 231.598946
$ws.Range("I25").Value = 0.2284507267316773
$ws.Range("J25").Value = 0.2432490894289675
$ws.Range("M25").Value = 0.1125055
$ws.Range("N25").Value = 0.225011
$ws.Range("O25").Value = 0.3657710320253459
$ws.Range("P25").Value = 0.2777063318885082
$ws.Range("Q25").Value = 8.685385073067666
$ws.Range("R25").Value = 52.112310438406
$ws.Range("S25").Value = 0.08356065808358588
$ws.Range("T25").Value = 0.06755181236053825
$ws.Range("G26").Value = 77.19964866666666
$ws.Range("H26").Value = 231.598946
$ws.Range("I26").Value = 0.2284507267316773
$ws.Range("J26").Value = 0.2432490894289675
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 0.3333333333333333
$ws.Range("M26").Value = 0.03062
$ws.Range("N26").Value = 0.09186
$ws.Range("O26").Value = 0.09954987978913112
$ws.Range("P26").Value = 0.1133726957672219
$ws.Range("Q26").Value = 2.363853242173333
$ws.Range("R26").Value = 21.27467917956
$ws.Range("S26").Value = 0.02274224238387811
$ws.Range("T26").Value = 0.01582573164992516
